$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C19").Value = "Henry’s Law constant (air/water partition coefficient) at 25 C"
$ws.Range("C20").Value = "Octanol-air partition coefficient"
$ws.Range("D22").Value = "degree C"
$ws.Range("C24").Value = "Water solubility at 25 C"
$ws.Range("C27").Value = "Biodegradation half-life for compounds containing only carbon and hydrogen "
$ws.Range("D28").Value = "binary (0/1)"
$ws.Range("C29").Value = "The whole body primary biotransformation rate (half-life) constant for organic chemicals in fish"
$ws.Range("C30").Value = "Soil adsorption coefficient of organic compounds"
$ws.Range("D31").Value = "binary (0/1)"
$ws.Range("D32").Value = "binary (0/2)"
$ws.Range("D33").Value = "binary (0/3)"
$ws.Range("D34").Value = "binary (0/4)"
$ws.Range("D35").Value = "binary (0/5)"
$ws.Range("D36").Value = "binary (0/6)"

$ws.Range("C2").Select()
